$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend column Y with the new 2021 data, copying the number/cell format
# from the corresponding cell in column X (2020) for every row, then
# filling in the actual 2021 values.
$ws.Range("X4:X16").Copy()
$ws.Range("Y4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("Y4").Value = 2021
$ws.Range("Y5").Value = 46.69
$ws.Range("Y6").Value = 52.52
$ws.Range("Y7").Value = 43.22
$ws.Range("Y8").Value = 51.31
$ws.Range("Y9").Value = 41.31
$ws.Range("Y10").Value = 52.43
$ws.Range("Y11").Value = 49.27
$ws.Range("Y12").Value = 31.68
$ws.Range("Y13").Value = 35.59
$ws.Range("Y14").Value = 55.28
$ws.Range("Y15").Value = 61.02
$ws.Range("Y16").Value = 48.72

# Restore the current selection to match the state captured when the
# workbook was last saved.
$ws.Range("AA15").Select()
